# feat: add 2022-Q1 data
#
# Before: sheets = [ "2021-Q4" (fund holdings), "总计" (summary) ]
# After:  sheets = [ "2021-Q4" (fund holdings, unchanged),
#                     "2022-Q1" (new fund holdings sheet),
#                     "总计" (summary, gains a new leading row for 2022-Q1) ]

function Set-TextValue($range, [string]$value) {
    # Force a numeric-looking string (e.g. "002802", "0.41") to be stored as
    # TEXT rather than auto-coerced to a number (which would also eat
    # leading zeros). The leading apostrophe is Excel's classic
    # text-prefix; resetting .Style back to "Normal" afterwards drops the
    # "number stored as text" quote-prefix flag/style that the apostrophe
    # trick leaves behind, so the cell ends up plain text with the default
    # style - matching a cell that was simply authored as text to begin with.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet by duplicating "2021-Q4"'s layout/
#        formatting, positioned right after it, then overwrite its values. ---
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

Set-TextValue $q1.Range("B2") "002802"
$q1.Range("C2").Value = "广发东财大数据精选灵活配置混合"
Set-TextValue $q1.Range("D2") "0.41"
Set-TextValue $q1.Range("E2") "55.13"
Set-TextValue $q1.Range("F2") "2.84"
Set-TextValue $q1.Range("G2") "0.0116"
$q1.Range("H2").Value = 2

# Re-fetch "总计" AFTER the sheet insert above: worksheet handles obtained
# by name/index before a Worksheets.Add/Copy can end up pointing at the
# newly-inserted sheet once indices shift, so look it up fresh here.
$total = $wb.Worksheets.Item("总计")

# --- 2. Update the "总计" summary sheet: write the existing 2021-Q4 row
#        into row 3 first, then overwrite row 2 with the new 2022-Q1 row
#        (avoids Rows.Insert(), which drags unwanted formatting along). ---
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.05

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

# Row 3's "A" cell is brand new (the old sheet only ever had one data row)
# so it starts out with the default style; clone A2's style (border/bold,
# matching the rest of column A) onto it without touching its value.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$total.Application.CutCopyMode = $false
